# Apply cryptos.xlsx price/volume update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D holding plain numeric-looking text need an explicit Text format
# first, otherwise Excel auto-converts the assigned string into a real number
# (losing the trailing zero / exact text representation from the source feed).
$textCells = @("D5","D6","D17","D19","D20","D21","D22","D24","D26","D29","D31","D32","D35","D36","D38","D39","D46","D47","D49","D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '72.358.19'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.656.52'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '592.57'
$ws.Range('E5').Value = '  -1.73%  '
$ws.Range('D6').Value = '174.30'
$ws.Range('E6').Value = '  -3.23%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('D9').Value = '2.656.25'
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('E10').Value = '  -2.58%  '
$ws.Range('E11').Value = '  +1.41%  '
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('E13').Value = '  -1.64%  '
$ws.Range('D14').Value = '3.143.28'
$ws.Range('E14').Value = '  +1.01%  '
$ws.Range('E15').Value = '  -2.37%  '
$ws.Range('D16').Value = '72.260.81'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').Value = '26.08'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').Value = '2.653.22'
$ws.Range('D19').Value = '12.39'
$ws.Range('E19').Value = '  +3.42%  '
$ws.Range('D20').Value = '8.02'
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('D21').Value = '372.52'
$ws.Range('E21').Value = '  -1.89%  '
$ws.Range('D22').Value = '4.18'
$ws.Range('E23').Value = '  -1.35%  '
$ws.Range('D24').Value = '71.69'
$ws.Range('E24').Value = '  -2.28%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').Value = '4.26'
$ws.Range('E26').Value = '  -3.19%  '
$ws.Range('E27').Value = '  -4.37%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').Value = '8.10'
$ws.Range('E31').Value = '  -0.62%  '
$ws.Range('D32').Value = '498.85'
$ws.Range('E32').Value = '  -4.91%  '
$ws.Range('E33').Value = '  -3.00%  '
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Value = '161.80'
$ws.Range('E36').Value = '  -1.52%  '
$ws.Range('E37').Value = '  +3.94%  '
$ws.Range('D38').Value = '19.42'
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('D39').Value = '18.90'
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('E40').Value = '  -3.44%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  -6.87%  '
$ws.Range('E43').Value = '  -3.87%  '
$ws.Range('E44').Value = '  -3.87%  '
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '39.18'
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '154.92'
$ws.Range('E47').Value = '  +2.19%  '
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('D49').Value = '0.551'
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.0749'
$ws.Range('E51').Value = '  -1.37%  '
